$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) look like numbers ("1.00", "0.0000240", "2.984.72", ...)
# and the Value setter would silently coerce them to a numeric type,
# dropping trailing zeros / exponent-formatting the text. Force Text format
# first so they round-trip as the literal strings from the source data.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D12","D13","D14","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D27","D29","D31","D32","D34","D36","D37","D38","D41","D43","D44","D45","D46","D47","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.967.98"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").Value = "2.965.04"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "595.79"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "146.50"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "2.960.19"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").Value = "0.508"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").Value = "7.23"
$ws.Range("E10").Value = "  +3.23%  "
$ws.Range("E11").Value = "  +6.50%  "
$ws.Range("D12").Value = "0.445"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "0.0000240"
$ws.Range("E13").Value = "  +6.63%  "
$ws.Range("D14").Value = "33.30"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "3.460.50"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").Value = "62.883.12"
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.984.72"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "6.74"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "442.14"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").Value = "13.58"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "0.673"
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("D23").Value = "7.09"
$ws.Range("D24").Value = "11.28"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("D25").Value = "81.62"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("D27").Value = "11.88"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "7.27"
$ws.Range("E29").Value = "  +3.98%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "2.16"
$ws.Range("E31").Value = "  -3.44%  "
$ws.Range("D32").Value = "0.0₃0970"
$ws.Range("E32").Value = "  +10.69%  "
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").Value = "26.58"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "0.994"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "3.13"
$ws.Range("E37").Value = "  +5.24%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "5.66"
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").Value = "8.53"
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("E42").Value = "  -3.75%  "
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D43").Value = "41.02"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.282"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "2.744.80"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").Value = "134.77"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").Value = "366.85"
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("D50").Value = "23.03"
$ws.Range("E50").Value = "  -3.67%  "
$ws.Range("E51").Value = "  -0.39%  "
